$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.805.14"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").Value = "1.873.97"
$ws.Range("E3").Value = "  +1.47%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9980"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4934"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.90"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2901"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06601"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.02%  "

$ws.Range("D11").Value = "1.876.09"
$ws.Range("E11").Value = "  +1.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.88"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07160"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("E14").Value = "  +1.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.32"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.798"
$ws.Range("D16").ClearFormats()

$ws.Range("D17").Value = "29.797.70"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007797"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +5.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9979"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.71"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.01%  "

$ws.Range("D21").Value = "2.119.60"
$ws.Range("E21").Value = "  +2.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9975"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.724"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.124"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.556"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.10"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "133.45"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.67"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.920"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.379"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.34%  "

$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08650"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.901"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05044"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7057"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.65%  "

$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.665"
$ws.Range("D37").ClearFormats()

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.202"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.655"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9289"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.08%  "

$ws.Range("E41").Value = "  +1.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.068"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9943"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.36"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4166"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.493"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1253"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05691"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.50"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.221"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.335"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.17%  "
